$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.987.94"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "3.402.85"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "581.95"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "138.86"
$ws.Range("E6").Value = "  +4.20%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.400.96"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").Value = "0.475"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "7.52"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("E11").Value = "  +5.33%  "
$ws.Range("D12").Value = "0.393"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").Value = "3.987.14"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").Value = "3.405.25"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "25.52"
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").Value = "62.087.72"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("D19").Value = "14.24"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("D20").Value = "9.54"
$ws.Range("E20").Value = "  +2.83%  "
$ws.Range("D21").Value = "5.83"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "394.37"
$ws.Range("E22").Value = "  +4.87%  "
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "0.0000131"
$ws.Range("E24").Value = "  +13.13%  "
$ws.Range("D25").Value = "3.550.95"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "71.68"
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("D28").Value = "1.66"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "7.70"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +4.37%  "
$ws.Range("D32").Value = "8.26"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "23.56"
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D36").Value = "3.437.29"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("D37").Value = "5.43"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("E38").Value = "  +3.14%  "
$ws.Range("D39").Value = "6.95"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "0.0791"
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("B42").Value = "ONDO"
$ws.Range("C42").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D42").Value = "1.28"
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.79"
$ws.Range("E43").Value = "  +11.51%  "
$ws.Range("D44").Value = "0.789"
$ws.Range("E44").Value = "  +4.99%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "4.46"
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "41.72"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "25.00"
$ws.Range("E48").Value = "  +6.39%  "
$ws.Range("D49").Value = "6.93"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "23.33"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("D51").Value = "2.351.24"
$ws.Range("E51").Value = "  +7.72%  "
